$wb = $excel.ActiveWorkbook

# The sheet formerly called "5_Quant" loses its "_Quant" suffix.
$ws = $wb.Worksheets.Item("5_Quant")
$ws.Name = "5_"

# Update the question text (A1), keep the Leeway/Comments headers (B1/C1) as-is.
$ws.Range("A1").Value = "What is a good estimate for the birth rate?"

# Update the example answer/leeway numbers in row 2.
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = 1

# Replace the single remaining comment (row 3) with the new guidance text.
$ws.Range("C3").Value = "You want the gray ""model"" line to fit the known data as closely as possible.  At this point, just ""eyeball"" a good approximation."

# The two extra comment rows (4 and 5) are no longer needed.
$ws.Range("A4:C5").EntireRow.Delete()

# Match the row heights used for the wrapped question/comment text.
$ws.Range("A1:C1").RowHeight = 30
$ws.Range("A3:C3").RowHeight = 60

# This sheet becomes the active tab, with B8 selected.
$ws.Activate() | Out-Null
$ws.Range("B8").Select() | Out-Null
